$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 45107
$ws.Range("D4").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104003
$ws.Range("J4").Value = "Membrillo"
$ws.Range("K4").Value = "Champion"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = '$/bandeja 18 kilos'
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1083
$ws.Range("T4").Value = 18
